$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "RunMachine"
$ws.Cells.Item(2, 2).Value = "Barbados Royals"
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 8
$ws.Cells.Item(2, 5).Value = 0.39
$ws.Cells.Item(2, 6).Value = 16
$ws.Cells.Item(3, 1).Value = "QuantumQuirk"
$ws.Cells.Item(3, 2).Value = "Peshawar Zalmi"
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = 7
$ws.Cells.Item(3, 5).Value = -0.04
$ws.Cells.Item(3, 6).Value = 14
$ws.Cells.Item(4, 1).Value = "GentleGamer"
$ws.Cells.Item(4, 2).Value = "Multan Sultans"
$ws.Cells.Item(4, 3).Value = 12
$ws.Cells.Item(4, 4).Value = 7
$ws.Cells.Item(4, 5).Value = -0.26
$ws.Cells.Item(4, 6).Value = 14
$ws.Cells.Item(5, 1).Value = "CodeCricketMaster"
$ws.Cells.Item(5, 2).Value = "Sydney Sixers"
$ws.Cells.Item(5, 3).Value = 12
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 0.6
$ws.Cells.Item(5, 6).Value = 12
$ws.Cells.Item(6, 1).Value = "db1_db2"
$ws.Cells.Item(6, 2).Value = "Chennai Super Kings"
$ws.Cells.Item(6, 3).Value = 12
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).Value = 0.32
$ws.Cells.Item(6, 6).Value = 12
$ws.Cells.Item(7, 1).Value = "Mike"
$ws.Cells.Item(7, 2).Value = "Jamaica Tallawahs"
$ws.Cells.Item(7, 3).Value = 12
$ws.Cells.Item(7, 4).Value = 5
$ws.Cells.Item(7, 5).Value = 0.74
$ws.Cells.Item(7, 6).Value = 10
$ws.Cells.Item(8, 1).Value = "Rahul"
$ws.Cells.Item(8, 2).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(8, 3).Value = 12
$ws.Cells.Item(8, 4).Value = 3
$ws.Cells.Item(8, 5).Value = -1.76
$ws.Cells.Item(8, 6).Value = 6
$ws.Cells.Item(9, 1).Value = "dosu"
$ws.Cells.Item(9, 2).Value = "Guyana Amazon Warriors"
$ws.Cells.Item(9, 3).Value = 12
$ws.Cells.Item(9, 4).Value = 8
$ws.Cells.Item(9, 5).Value = 0.22
$ws.Cells.Item(9, 6).Value = 16
$ws.Cells.Item(10, 1).Value = "Ved"
$ws.Cells.Item(10, 2).Value = "Kolkata Knight Riders"
$ws.Cells.Item(10, 3).Value = 12
$ws.Cells.Item(10, 4).Value = 7
$ws.Cells.Item(10, 5).Value = 0.45
$ws.Cells.Item(10, 6).Value = 14
$ws.Cells.Item(11, 1).Value = "sammat"
$ws.Cells.Item(11, 2).Value = "Karachi Kings"
$ws.Cells.Item(11, 3).Value = 12
$ws.Cells.Item(11, 4).Value = 7
$ws.Cells.Item(11, 5).Value = 0.42
$ws.Cells.Item(11, 6).Value = 14
$ws.Cells.Item(12, 1).Value = "newGuy"
$ws.Cells.Item(12, 2).Value = "St Lucia Kings"
$ws.Cells.Item(12, 3).Value = 12
$ws.Cells.Item(12, 4).Value = 5
$ws.Cells.Item(12, 5).Value = 0.27
$ws.Cells.Item(12, 6).Value = 10
$ws.Cells.Item(13, 1).Value = "Crabby"
$ws.Cells.Item(13, 2).Value = "Hobart Hurricanes"
$ws.Cells.Item(13, 3).Value = 12
$ws.Cells.Item(13, 4).Value = 5
$ws.Cells.Item(13, 5).Value = -0.1
$ws.Cells.Item(13, 6).Value = 10
$ws.Cells.Item(14, 1).Value = "Sricharan"
$ws.Cells.Item(14, 2).Value = "Melbourne Stars"
$ws.Cells.Item(14, 3).Value = 12
$ws.Cells.Item(14, 4).Value = 5
$ws.Cells.Item(14, 5).Value = -0.44
$ws.Cells.Item(14, 6).Value = 10
$ws.Cells.Item(15, 1).Value = "Sachin"
$ws.Cells.Item(15, 2).Value = "Melbourne Renegades"
$ws.Cells.Item(15, 3).Value = 12
$ws.Cells.Item(15, 4).Value = 5
$ws.Cells.Item(15, 5).Value = -0.86
$ws.Cells.Item(15, 6).Value = 10
$ws.Cells.Item(16, 1).Value = "Angel"
$ws.Cells.Item(16, 2).Value = "Delhi Capitals"
$ws.Cells.Item(16, 3).Value = 14
$ws.Cells.Item(16, 4).Value = 9
$ws.Cells.Item(16, 5).Value = 0.32
$ws.Cells.Item(16, 6).Value = 18
$ws.Cells.Item(17, 1).Value = "Ava"
$ws.Cells.Item(17, 2).Value = "Punjab Kings"
$ws.Cells.Item(17, 3).Value = 14
$ws.Cells.Item(17, 4).Value = 8
$ws.Cells.Item(17, 5).Value = 0.81
$ws.Cells.Item(17, 6).Value = 16
$ws.Cells.Item(18, 1).Value = "ady_chak"
$ws.Cells.Item(18, 2).Value = "Mumbai Indians"
$ws.Cells.Item(18, 3).Value = 14
$ws.Cells.Item(18, 4).Value = 8
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 16
$ws.Cells.Item(19, 1).Value = "Maverick"
$ws.Cells.Item(19, 2).Value = "Brisbane Heat"
$ws.Cells.Item(19, 3).Value = 14
$ws.Cells.Item(19, 4).Value = 7
$ws.Cells.Item(19, 5).Value = 0.35
$ws.Cells.Item(19, 6).Value = 14
$ws.Cells.Item(20, 1).Value = "Prateesh"
$ws.Cells.Item(20, 2).Value = "Perth Scorchers"
$ws.Cells.Item(20, 3).Value = 14
$ws.Cells.Item(20, 4).Value = 6
$ws.Cells.Item(20, 5).Value = 0.37
$ws.Cells.Item(20, 6).Value = 12
$ws.Cells.Item(21, 1).Value = "GeniiExE"
$ws.Cells.Item(21, 2).Value = "Adelaide Strikers"
$ws.Cells.Item(21, 3).Value = 14
$ws.Cells.Item(21, 4).Value = 6
$ws.Cells.Item(21, 5).Value = 0.13
$ws.Cells.Item(21, 6).Value = 12
$ws.Cells.Item(22, 1).Value = "AnkitGamer"
$ws.Cells.Item(22, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(22, 3).Value = 14
$ws.Cells.Item(22, 4).Value = 6
$ws.Cells.Item(22, 5).Value = -0.98
$ws.Cells.Item(22, 6).Value = 12
$ws.Cells.Item(23, 1).Value = "Bawandar"
$ws.Cells.Item(23, 2).Value = "Sydney Thunder"
$ws.Cells.Item(23, 3).Value = 14
$ws.Cells.Item(23, 4).Value = 6
$ws.Cells.Item(23, 5).Value = -1.02
$ws.Cells.Item(23, 6).Value = 12

$ws.Range("J23").Select()
